# Scheduled runner update: refresh market-price-derived leve profit figures
# (currentAveragePrice*, LevePriceNQ/HQ, LeveProfitNQ/HQ columns) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1065.875
$ws.Range("I43").Value = 1064.1666
$ws.Range("K43").Value = 1064.1666
$ws.Range("M43").Value = -995.1666
$ws.Range("H76").Value = 6007
$ws.Range("I76").Value = 6115
$ws.Range("J76").Value = 5899
$ws.Range("K76").Value = 6115
$ws.Range("L76").Value = 5899
$ws.Range("M76").Value = -5800
$ws.Range("N76").Value = -6529
$ws.Range("H79").Value = 6007
$ws.Range("I79").Value = 6115
$ws.Range("J79").Value = 5899
$ws.Range("K79").Value = 6115
$ws.Range("L79").Value = 5899
$ws.Range("M79").Value = -5023
$ws.Range("N79").Value = -8083
$ws.Range("H80").Value = 845.9666999999999
$ws.Range("J80").Value = 1287.7142
$ws.Range("L80").Value = 3863.1426
$ws.Range("N80").Value = -5859.142599999999
$ws.Range("H83").Value = 845.9666999999999
$ws.Range("J83").Value = 1287.7142
$ws.Range("L83").Value = 11589.4278
$ws.Range("N83").Value = -21573.4278
$ws.Range("H100").Value = 566.6667
$ws.Range("I100").Value = 600
$ws.Range("J100").Value = 550
$ws.Range("K100").Value = 600
$ws.Range("L100").Value = 550
$ws.Range("M100").Value = -59
$ws.Range("N100").Value = -1632
$ws.Range("H137").Value = 120788.53
$ws.Range("I137").Value = 179575.8
$ws.Range("J137").Value = 3214
$ws.Range("K137").Value = 538727.3999999999
$ws.Range("L137").Value = 9642
$ws.Range("M137").Value = -536177.3999999999
$ws.Range("N137").Value = -14742

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5477.577
$ws.Range("I32").Value = 3187.9
$ws.Range("K32").Value = 3187.9
$ws.Range("M32").Value = -2900.9
$ws.Range("H74").Value = 35133.723
$ws.Range("I74").Value = 4991.66
$ws.Range("J74").Value = 172143.1
$ws.Range("K74").Value = 4991.66
$ws.Range("L74").Value = 172143.1
$ws.Range("M74").Value = -4117.66
$ws.Range("N74").Value = -173891.1
$ws.Range("H77").Value = 35133.723
$ws.Range("I77").Value = 4991.66
$ws.Range("J77").Value = 172143.1
$ws.Range("K77").Value = 24958.3
$ws.Range("L77").Value = 860715.5
$ws.Range("M77").Value = -20590.3
$ws.Range("N77").Value = -869451.5
$ws.Range("H102").Value = 4633010
$ws.Range("I102").Value = 5954918.5
$ws.Range("K102").Value = 5954918.5
$ws.Range("M102").Value = -5953296.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7146954
$ws.Range("I99").Value = 11908349
$ws.Range("K99").Value = 11908349
$ws.Range("M99").Value = -11906851
$ws.Range("H105").Value = 4167996
$ws.Range("I105").Value = 4465653
$ws.Range("J105").Value = 797
$ws.Range("K105").Value = 4465653
$ws.Range("L105").Value = 797
$ws.Range("M105").Value = -4463906
$ws.Range("N105").Value = -4291

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4376.125
$ws.Range("J58").Value = 4502.3335
$ws.Range("L58").Value = 4502.3335
$ws.Range("N58").Value = -4908.3335
$ws.Range("H105").Value = 2030
$ws.Range("I105").Value = 1795
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1795
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -48
$ws.Range("N105").Value = -5994
$ws.Range("H136").Value = 4376.125
$ws.Range("J136").Value = 4502.3335
$ws.Range("L136").Value = 13507.0005
$ws.Range("N136").Value = -18607.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4472831
$ws.Range("I4").Value = 6964603.5
$ws.Range("K4").Value = 20893810.5
$ws.Range("M4").Value = -20893698.5
$ws.Range("H14").Value = 1209.2
$ws.Range("I14").Value = 1209.2
$ws.Range("K14").Value = 3627.6
$ws.Range("M14").Value = -3454.6
$ws.Range("H88").Value = 14999.833
$ws.Range("J88").Value = 14999.833
$ws.Range("L88").Value = 44999.499
$ws.Range("N88").Value = -45855.499
$ws.Range("H91").Value = 14999.833
$ws.Range("J91").Value = 14999.833
$ws.Range("L91").Value = 44999.499
$ws.Range("N91").Value = -47963.499
$ws.Range("H118").Value = 2302.8
$ws.Range("I118").Value = 2066
$ws.Range("J118").Value = 3250
$ws.Range("K118").Value = 6198
$ws.Range("L118").Value = 9750
$ws.Range("M118").Value = -4955
$ws.Range("N118").Value = -12236
$ws.Range("H121").Value = 5829.091
$ws.Range("I121").Value = 7702.5
$ws.Range("J121").Value = 833.3333
$ws.Range("K121").Value = 23107.5
$ws.Range("L121").Value = 2499.9999
$ws.Range("M121").Value = -21797.5
$ws.Range("N121").Value = -5119.9999
$ws.Range("H132").Value = 1982.92
$ws.Range("I132").Value = 1220.9286
$ws.Range("J132").Value = 2952.7273
$ws.Range("K132").Value = 10988.3574
$ws.Range("L132").Value = 26574.5457
$ws.Range("M132").Value = -8458.357399999999
$ws.Range("N132").Value = -31634.5457
$ws.Range("H136").Value = 1646.5714
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1358336.8
$ws.Range("I80").Value = 2441275
$ws.Range("J80").Value = 4663.875
$ws.Range("K80").Value = 2441275
$ws.Range("L80").Value = 4663.875
$ws.Range("M80").Value = -2440277
$ws.Range("N80").Value = -6659.875
$ws.Range("H83").Value = 1358336.8
$ws.Range("I83").Value = 2441275
$ws.Range("J83").Value = 4663.875
$ws.Range("K83").Value = 12206375
$ws.Range("L83").Value = 23319.375
$ws.Range("M83").Value = -12201383
$ws.Range("N83").Value = -33303.375
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744
$ws.Range("H132").Value = 2572.04
$ws.Range("I132").Value = 1927.3
$ws.Range("J132").Value = 5151
$ws.Range("K132").Value = 5781.9
$ws.Range("L132").Value = 15453
$ws.Range("M132").Value = -3251.9
$ws.Range("N132").Value = -20513

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 462.26666
$ws.Range("I16").Value = 488.25925
$ws.Range("J16").Value = 228.33333
$ws.Range("K16").Value = 488.25925
$ws.Range("L16").Value = 228.33333
$ws.Range("M16").Value = -318.25925
$ws.Range("N16").Value = -568.3333299999999
$ws.Range("H46").Value = 4164.2354
$ws.Range("I46").Value = 899
$ws.Range("J46").Value = 6449.9
$ws.Range("K46").Value = 899
$ws.Range("L46").Value = 6449.9
$ws.Range("M46").Value = -711
$ws.Range("N46").Value = -6825.9
$ws.Range("H100").Value = 3500
$ws.Range("J100").Value = 8000
$ws.Range("L100").Value = 8000
$ws.Range("N100").Value = -9082
$ws.Range("H136").Value = 66805.69
$ws.Range("I136").Value = 94717.55
$ws.Range("J136").Value = 5399.6
$ws.Range("K136").Value = 284152.65
$ws.Range("L136").Value = 16198.8
$ws.Range("M136").Value = -281602.65
$ws.Range("N136").Value = -21298.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H52").Value = 13180.667
$ws.Range("I52").Value = 9771
$ws.Range("J52").Value = 20000
$ws.Range("K52").Value = 9771
$ws.Range("L52").Value = 20000
$ws.Range("M52").Value = -9545
$ws.Range("N52").Value = -20452
$ws.Range("H81").Value = 55556624
$ws.Range("J81").Value = 1600
$ws.Range("L81").Value = 3200
$ws.Range("N81").Value = -5322
$ws.Range("H84").Value = 55556624
$ws.Range("J84").Value = 1600
$ws.Range("L84").Value = 16000
$ws.Range("N84").Value = -26608
$ws.Range("H96").Value = 1000
$ws.Range("J96").Value = 1000
$ws.Range("L96").Value = 1000
$ws.Range("N96").Value = -3746
